$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = 3914
$ws.Range("B4").Value = 82
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4779

$ws.Range("A8").Value = 15721
$ws.Range("B8").Value = 301
$ws.Range("A9").Value = 16
$ws.Range("B9").Value = 19078

$ws.Range("A14").Value = "         0.0       1.00      0.98      0.99      3996"
$ws.Range("A15").Value = "         1.0       0.98      1.00      0.99      4783"
$ws.Range("A17").Value = "    accuracy                           0.99      8779"
$ws.Range("A18").Value = "   macro avg       0.99      0.99      0.99      8779"
$ws.Range("A19").Value = "weighted avg       0.99      0.99      0.99      8779"

$ws.Range("A25").Value = "         0.0       1.00      0.98      0.99     16022"
$ws.Range("A26").Value = "         1.0       0.98      1.00      0.99     19094"
$ws.Range("A28").Value = "    accuracy                           0.99     35116"
$ws.Range("A29").Value = "   macro avg       0.99      0.99      0.99     35116"
$ws.Range("A30").Value = "weighted avg       0.99      0.99      0.99     35116"

$ws.Range("B33").Value = 0.9902038956600979
$ws.Range("B35").Value = 0.9909727759425903
$ws.Range("B37").Value = 0.9831310429952684
$ws.Range("B39").Value = 0.9844677227927138
$ws.Range("B41").Value = 0.99916370478779
$ws.Range("B43").Value = 0.9991620404315492
$ws.Range("B45").Value = 0.9910825383658233
$ws.Range("B47").Value = 0.9917604553842955

$ws.Range("B49").Value = "0:02:46.388099"
$ws.Range("B51").Value = "{'solver': 'adam', 'hidden_layer_sizes': (100,), 'alpha': 0.001, 'activation': 'logistic'}"
